# Auto update stock data: bump the "Date_1" column from 2025/12/25 to
# 2025/12/26 for every company block's first row.
#
# The cells hold the date as plain text (inline string), not a real Excel
# date serial. Assigning a date-shaped string straight to .Value makes the
# engine auto-convert it into a date serial (and stamps a number-format
# style on the cell), which would change both the stored value type and
# the cell's style - neither of which happened in the source edit. To
# avoid that, we briefly force the cell to Text format before writing the
# literal string, then clear the (now unneeded) formatting so the cell's
# style index goes right back to the default, exactly like the original.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rows = @(2, 8, 14, 20, 26, 32, 38, 44, 50, 56, 62, 68, 74)

foreach ($r in $rows) {
    $cell = $ws.Cells.Item($r, 1)
    $cell.NumberFormat = "@"
    $cell.Value = "2025/12/26"
    $cell.ClearFormats()
}
